$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "859×9=7731" "229×2=458"
Replace-Text "283×4=1132" "981×7=6867"
Replace-Text "359×7=2513" "730×6=4380"
Replace-Text "604×8=4832" "588×9=5292"
Replace-Text "175×4=700" "324×2=648"
Replace-Text "430×2=860" "815×5=4075"
Replace-Text "268×9=2412" "920×2=1840"
Replace-Text "819×7=5733" "867×2=1734"
Replace-Text "193×3=579" "909×2=1818"
Replace-Text "938×6=5628" "765×7=5355"
Replace-Text "655×7=4585" "654×2=1308"
Replace-Text "912×3=2736" "515×3=1545"
Replace-Text "662×5=3310" "238×7=1666"
Replace-Text "688×8=5504" "450×2=900"
Replace-Text "285×2=570" "805×3=2415"
Replace-Text "791×4=3164" "277×5=1385"
Replace-Text "930×7=6510" "291×3=873"
Replace-Text "644×5=3220" "651×3=1953"
Replace-Text "475×3=1425" "116×7=812"
Replace-Text "761×7=5327" "472×2=944"
Replace-Text "767×9=6903" "628×6=3768"
Replace-Text "299×5=1495" "496×9=4464"
Replace-Text "683×5=3415" "227×7=1589"
Replace-Text "680×4=2720" "582×6=3492"
Replace-Text "376×4=1504" "421×3=1263"
